# Hortaliza, Mapocho Venta Directa de Santiago - Pepino dulce
# Weekly update: insert 3 new rows of data (most recent date 44460) at the
# top of the data block that starts at row 82, pushing the existing rows
# 82-129 down to 85-132. All other existing data is unchanged; it simply
# moves down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the old row 82 (this shifts rows 82:129 -> 85:132,
# and also grows the sheet's used range / dimension to A1:R132 automatically).
$ws.Range("A82:A84").EntireRow.Insert()

# Seed the 3 new rows with the same static/boilerplate column formatting and
# values as the rest of the block (columns A,B,C,E,F,G,H,O,Q,R are constant
# across every data row), by copying the row right below them (the row that
# used to be row 82, now shifted to row 85) down into each new row in turn
# (PasteSpecial only fills the first row of a multi-row destination, so
# paste once per row).
$ws.Range("A85:R85").Copy()
$ws.Range("A82:R82").PasteSpecial()
$ws.Range("A85:R85").Copy()
$ws.Range("A83:R83").PasteSpecial()
$ws.Range("A85:R85").Copy()
$ws.Range("A84:R84").PasteSpecial()

# --- Row 82: date 44460, Especial ---
$ws.Range("D82").Value2 = 44460
$ws.Range("I82").Value = "Especial"
$ws.Range("J82").Value2 = 450
$ws.Range("K82").Value2 = 22000
$ws.Range("L82").Value2 = 22000
$ws.Range("M82").Value2 = 22000
$ws.Range("N82").Value = "$/bandeja 18 kilos"
$ws.Range("P82").Value2 = 1222

# --- Row 83: date 44460, Primera ---
$ws.Range("D83").Value2 = 44460
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value2 = 510
$ws.Range("K83").Value2 = 20000
$ws.Range("L83").Value2 = 20000
$ws.Range("M83").Value2 = 20000
$ws.Range("N83").Value = "$/bandeja 18 kilos"
$ws.Range("P83").Value2 = 1111

# --- Row 84: date 44460, Segunda ---
$ws.Range("D84").Value2 = 44460
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value2 = 400
$ws.Range("K84").Value2 = 17000
$ws.Range("L84").Value2 = 17000
$ws.Range("M84").Value2 = 17000
$ws.Range("N84").Value = "$/bandeja 18 kilos"
$ws.Range("P84").Value2 = 944

